# Generate Report for Handoff
#
# The localization-status report moved from "In Translation" to
# "Ready for handoff", and the handoff timestamps were refreshed.
# This touches the three worksheets (Overview, zh-cn, de-de): the
# Status text, the two related "Handoff Datetime" timestamps, and the
# (now wider, to fit the new text) Status/zh-cn/de-de columns.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Status: "In Translation" -> "Ready for handoff" -----------------
$wsOverview.Range("E2").Value = "Ready for handoff"   # zh-cn status
$wsOverview.Range("F2").Value = "Ready for handoff"   # de-de status
$wsZhCn.Range("C2").Value     = "Ready for handoff"
$wsDeDe.Range("C2").Value     = "Ready for handoff"

# --- Refreshed "Latest Handoff Datetime" timestamps -------------------
$wsZhCn.Range("H2").Value     = "2016-08-17 00:54:41"
$wsDeDe.Range("H2").Value     = "2016-08-17 00:54:46"
$wsOverview.Range("G2").Value = "2016-08-17 00:54:46"

# --- Widen the Status columns so the longer text fits -----------------
$wsOverview.Columns.Item(5).ColumnWidth = 16.38
$wsOverview.Columns.Item(6).ColumnWidth = 16.38
$wsZhCn.Columns.Item(3).ColumnWidth     = 16.38
$wsDeDe.Columns.Item(3).ColumnWidth     = 16.38
